$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1154.8334
$ws.Range("I40").Value = 795.5
$ws.Range("J40").Value = 1657.9
$ws.Range("K40").Value = 795.5
$ws.Range("L40").Value = 1657.9
$ws.Range("M40").Value = -620.5
$ws.Range("N40").Value = -2007.9

$ws.Range("H55").Value = 73.41936
$ws.Range("I55").Value = 39.42857
$ws.Range("J55").Value = 83.333336
$ws.Range("K55").Value = 39.42857
$ws.Range("L55").Value = 83.333336
$ws.Range("M55").Value = 174.57143
$ws.Range("N55").Value = -511.333336

$ws.Range("H92").Value = 913.7778
$ws.Range("I92").Value = 1049.1428
$ws.Range("K92").Value = 1049.1428
$ws.Range("M92").Value = 198.8571999999999

$ws.Range("H99").Value = 225.81818
$ws.Range("I99").Value = 197.71428
$ws.Range("J99").Value = 275
$ws.Range("K99").Value = 593.14284
$ws.Range("L99").Value = 825
$ws.Range("M99").Value = 904.85716
$ws.Range("N99").Value = -3821

$ws.Range("H112").Value = 1125.9
$ws.Range("J112").Value = 1149.25
$ws.Range("L112").Value = 3447.75
$ws.Range("N112").Value = -5663.75

$ws.Range("H116").Value = 10419754
$ws.Range("I116").Value = 20834690
$ws.Range("J116").Value = 4817.6665
$ws.Range("K116").Value = 20834690
$ws.Range("L116").Value = 4817.6665
$ws.Range("M116").Value = -20831248
$ws.Range("N116").Value = -11701.6665

$ws.Range("H137").Value = 81710.38
$ws.Range("I137").Value = 112275.39
$ws.Range("J137").Value = 3114.6428
$ws.Range("K137").Value = 336826.17
$ws.Range("L137").Value = 9343.928400000001
$ws.Range("M137").Value = -334276.17
$ws.Range("N137").Value = -14443.9284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H32").Value = 8001.8184
$ws.Range("I32").Value = 5806.4136
$ws.Range("J32").Value = 20667.615
$ws.Range("K32").Value = 5806.4136
$ws.Range("L32").Value = 20667.615
$ws.Range("M32").Value = -5519.4136
$ws.Range("N32").Value = -21241.615

$ws.Range("H45").Value = 2302.0938
$ws.Range("I45").Value = 2039.95
$ws.Range("J45").Value = 2739
$ws.Range("K45").Value = 2039.95
$ws.Range("L45").Value = 2739
$ws.Range("M45").Value = -1662.95
$ws.Range("N45").Value = -3493

$ws.Range("H61").Value = 2166.9
$ws.Range("I61").Value = 1818.8148
$ws.Range("J61").Value = 5299.6665
$ws.Range("K61").Value = 1818.8148
$ws.Range("L61").Value = 5299.6665
$ws.Range("M61").Value = -1606.8148
$ws.Range("N61").Value = -5723.6665

$ws.Range("H88").Value = 251040.25
$ws.Range("J88").Value = 334053.66
$ws.Range("L88").Value = 334053.66
$ws.Range("N88").Value = -334865.66

$ws.Range("H91").Value = 251040.25
$ws.Range("J91").Value = 334053.66
$ws.Range("L91").Value = 334053.66
$ws.Range("N91").Value = -336861.66

$ws.Range("H136").Value = 2166.9
$ws.Range("I136").Value = 1818.8148
$ws.Range("J136").Value = 5299.6665
$ws.Range("K136").Value = 5456.4444
$ws.Range("L136").Value = 15898.9995
$ws.Range("M136").Value = -2906.4444
$ws.Range("N136").Value = -20998.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 22939.8
$ws.Range("I96").Value = 5864
$ws.Range("J96").Value = 34323.668
$ws.Range("K96").Value = 5864
$ws.Range("L96").Value = 34323.668
$ws.Range("M96").Value = -3118
$ws.Range("N96").Value = -39815.668

$ws.Range("H99").Value = 1498.5
$ws.Range("I99").Value = 1397.1428
$ws.Range("J99").Value = 1640.4
$ws.Range("K99").Value = 1397.1428
$ws.Range("L99").Value = 1640.4
$ws.Range("M99").Value = 100.8571999999999
$ws.Range("N99").Value = -4636.4

$ws.Range("H134").Value = 3539.524
$ws.Range("I134").Value = 3504.3901
$ws.Range("J134").Value = 4980
$ws.Range("K134").Value = 10513.1703
$ws.Range("L134").Value = 14940
$ws.Range("M134").Value = -7978.1703
$ws.Range("N134").Value = -20010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3562.4736
$ws.Range("I31").Value = 1613.1936
$ws.Range("J31").Value = 5886.615
$ws.Range("K31").Value = 1613.1936
$ws.Range("L31").Value = 5886.615
$ws.Range("M31").Value = -1318.1936
$ws.Range("N31").Value = -6476.615

$ws.Range("H34").Value = 3562.4736
$ws.Range("I34").Value = 1613.1936
$ws.Range("J34").Value = 5886.615
$ws.Range("K34").Value = 1613.1936
$ws.Range("L34").Value = 5886.615
$ws.Range("M34").Value = -1411.1936
$ws.Range("N34").Value = -6290.615

$ws.Range("H62").Value = 4100
$ws.Range("J62").Value = 3750
$ws.Range("L62").Value = 3750
$ws.Range("N62").Value = -4998

$ws.Range("H65").Value = 4100
$ws.Range("J65").Value = 3750
$ws.Range("L65").Value = 18750
$ws.Range("N65").Value = -24990

$ws.Range("H86").Value = 33356346
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 41694684
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 41694684
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -41696930

$ws.Range("H89").Value = 33356346
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 41694684
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 208473420
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -208484652

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 810.1905
$ws.Range("I113").Value = 566.6667
$ws.Range("K113").Value = 1700.0001
$ws.Range("M113").Value = 469.9999

$ws.Range("H131").Value = 621.8763
$ws.Range("I131").Value = 305.10345
$ws.Range("J131").Value = 756.9706
$ws.Range("K131").Value = 915.31035
$ws.Range("L131").Value = 2270.9118
$ws.Range("M131").Value = 4124.68965
$ws.Range("N131").Value = -12350.9118

$ws.Range("H137").Value = 9807571
$ws.Range("J137").Value = 14497623
$ws.Range("L137").Value = 43492869
$ws.Range("N137").Value = -43503069

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2871.2856
$ws.Range("I7").Value = 2627.0908
$ws.Range("K7").Value = 2627.0908
$ws.Range("M7").Value = -2515.0908

$ws.Range("H126").Value = 2871.2856
$ws.Range("I126").Value = 2627.0908
$ws.Range("K126").Value = 7881.2724
$ws.Range("M126").Value = -5411.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1333.3334
$ws.Range("I122").Value = 1333.3334
$ws.Range("K122").Value = 4000.0002
$ws.Range("M122").Value = -1550.0002

Write-Host "edits applied"
